$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 103; this shifts the existing rows 103..178
# down to 104..179, matching the rest of the data already in the sheet.
$ws.Rows.Item(103).Insert()

# Populate the newly inserted row 103 with the new data record
# (same Mercado/Region/Categoria/etc. as its neighbours, new Fecha/Volumen).
$ws.Cells.Item(103, 1).Value = 10
$ws.Cells.Item(103, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(103, 3).Value = "La Araucanía"
$ws.Cells.Item(103, 4).Value = 44767
$ws.Cells.Item(103, 5).Value = 9
$ws.Cells.Item(103, 6).Value = 100114007
$ws.Cells.Item(103, 7).Value = "Jengibre"
$ws.Cells.Item(103, 8).Value = "Sin especificar"
$ws.Cells.Item(103, 9).Value = "Primera"
$ws.Cells.Item(103, 10).Value = 50
$ws.Cells.Item(103, 11).Value = 20000
$ws.Cells.Item(103, 12).Value = 20000
$ws.Cells.Item(103, 13).Value = 20000
$ws.Cells.Item(103, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(103, 15).Value = "Perú"
$ws.Cells.Item(103, 16).Value = 1538
$ws.Cells.Item(103, 17).Value = 13
$ws.Cells.Item(103, 18).Value = "Hortaliza"
